# Generate Report for handback
# Adds the handback status for e78d4cf2-7fdd-4f82-9714-0c81abb20e7c
# (source file + its zh-cn and de-de handback xlf) as a new row 4 on each
# of the three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$sourceMd   = "e78d4cf2-7fdd-4f82-9714-0c81abb20e7c.md"
$zhXlf      = "e78d4cf2-7fdd-4f82-9714-0c81abb20e7c.30251bb289d2b5e17899f9c1ed1dc147284e9a46.zh-cn.xlf"
$deXlf      = "e78d4cf2-7fdd-4f82-9714-0c81abb20e7c.30251bb289d2b5e17899f9c1ed1dc147284e9a46.de-de.xlf"
$statusSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $sourceMd
$wsOverview.Range("B4").Value = $statusSync
$wsOverview.Range("C4").Value = $statusSync

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/e78d4cf2-7fdd-4f82-9714-0c81abb20e7c.md",
    "",
    "",
    $sourceMd
)

# ---------------------------------------------------------------------
# zh-cn sheet: Source File Name | Status | Correspond Handoff File |
#   Correspond Handoff Datetime | Target File | Correspond Handback File |
#   Correspond Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $sourceMd
$wsZh.Range("B4").Value = $statusSync
$wsZh.Range("C4").Value = $zhXlf
$wsZh.Range("D4").Value = "2016-01-26 05:29:15"
$wsZh.Range("E4").Value = $sourceMd
$wsZh.Range("F4").Value = $zhXlf
$wsZh.Range("G4").Value = "2016-01-26 05:30:00"
$wsZh.Range("H4").Value = "Include"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/e78d4cf2-7fdd-4f82-9714-0c81abb20e7c.md",
    "",
    "",
    $sourceMd
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/e78d4cf2-7fdd-4f82-9714-0c81abb20e7c.30251bb289d2b5e17899f9c1ed1dc147284e9a46.zh-cn.xlf",
    "",
    "",
    $zhXlf
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/e78d4cf2-7fdd-4f82-9714-0c81abb20e7c.md",
    "",
    "",
    $sourceMd
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/e78d4cf2-7fdd-4f82-9714-0c81abb20e7c.30251bb289d2b5e17899f9c1ed1dc147284e9a46.zh-cn.xlf",
    "",
    "",
    $zhXlf
)

# ---------------------------------------------------------------------
# de-de sheet: same columns as zh-cn, but for the de-de handback file
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $sourceMd
$wsDe.Range("B4").Value = $statusSync
$wsDe.Range("C4").Value = $deXlf
$wsDe.Range("D4").Value = "2016-01-26 05:29:27"
$wsDe.Range("E4").Value = $sourceMd
$wsDe.Range("F4").Value = $deXlf
$wsDe.Range("G4").Value = "2016-01-26 05:30:21"
$wsDe.Range("H4").Value = "Include"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/e78d4cf2-7fdd-4f82-9714-0c81abb20e7c.md",
    "",
    "",
    $sourceMd
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/e78d4cf2-7fdd-4f82-9714-0c81abb20e7c.30251bb289d2b5e17899f9c1ed1dc147284e9a46.de-de.xlf",
    "",
    "",
    $deXlf
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/e78d4cf2-7fdd-4f82-9714-0c81abb20e7c.md",
    "",
    "",
    $sourceMd
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/e78d4cf2-7fdd-4f82-9714-0c81abb20e7c.30251bb289d2b5e17899f9c1ed1dc147284e9a46.de-de.xlf",
    "",
    "",
    $deXlf
)
